$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- B1: e000 content - IsChecked -> IsEnabled, split line before "Begin" button ---
$textB1 = @'
<Bold>e000 Welcome to Patton's Best Solo Tank Battle Game</Bold>
<LineBreak/><LineBreak/> <InlineUIContainer><CheckBox Name='SkipTutorial0' IsEnabled='True' Content='Skip in Future' FontFamily='Courier New'  FontSize='12'></CheckBox></InlineUIContainer>
<LineBreak/><LineBreak/>
The game starts with a tutorial how to play. However, before starting, it is important to know that Active events are shown with a green background. The game may only advance when a green background is displayed. Most often, the game advances by rolling dice or clicking an image. 
<LineBreak/><LineBreak/>
You can explore what may happen by showing inactive events. Inactive events have a gray background. To return to the current active event, select the active event button in the status bar per the image.
<LineBreak/>
               <InlineUIContainer><Image Name='Tutorial0' Height='70'  Width='370'> </Image></InlineUIContainer>
<LineBreak/><LineBreak/>
<InlineUIContainer><Button Name='Read_Rules' Content='Read Rules' FontFamily='Courier New'  FontSize='12'> </Button></InlineUIContainer> or 
<InlineUIContainer><Button Name='Begin' Content='Begin Game' FontFamily='Courier New'  FontSize='12'> </Button></InlineUIContainer>
'@

# --- B2: e001 content - IsChecked -> IsEnabled ---
$textB2 = @'
<Bold>e001 Fourth Armor Division Campaign</Bold> <InlineUIContainer><Button Content='r1.1' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>
<LineBreak/><LineBreak/> <InlineUIContainer><CheckBox Name='SkipTutorial1' IsEnabled='True' Content='Skip in Future' FontFamily='Courier New'  FontSize='12'></CheckBox></InlineUIContainer>
<LineBreak/><LineBreak/>
The campaign game of <Bold>Patton' Best</Bold> recreates the actions of the 4th Armored Division from late July 1944 through April 1945. 
<LineBreak/><LineBreak/>
Each day begins with a check of the Combat 
<InlineUIContainer><Button Content='Calendar' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> 
to see what the Division was doing on that day. The four possibilities are Refitting 
<InlineUIContainer><Button Content='r27.0' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>, 
an Advance scenario <InlineUIContainer><Button Content='r20.2' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>, a Battle scenario 
<InlineUIContainer><Button Content='r20.3' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>, or a Counterattack scenario 
<InlineUIContainer><Button Content='r20.4' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>.  Click image to continue.
<LineBreak/><LineBreak/>
                                            <InlineUIContainer><Image Name='Continue001' Height='100' Width='100'></Image></InlineUIContainer>
'@

# --- B3: e002 content - IsChecked -> IsEnabled ---
$textB3 = @'
<Bold>e002 Movement Board</Bold> <InlineUIContainer><Button Content='r2.11' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>
<LineBreak/><LineBreak/> <InlineUIContainer><CheckBox Name='SkipTutorial2' IsEnabled='True' Content='Skip in Future' FontFamily='Courier New'  FontSize='12'></CheckBox></InlineUIContainer>
<LineBreak/><LineBreak/>
The movement board is a depiction fo typical European countryside and is used to show the "big picture" for the day. The movement board is divided into white lines into areas. Click image to continue.
<LineBreak/><LineBreak/>
A=Farms    B=Fields    C=Villiages  D=Woods<LineBreak/>
#=Starting or exiting areas
<LineBreak/><LineBreak/>
                                   <InlineUIContainer><Image Name='MapMovement'  Height='200' Width='200'></Image></InlineUIContainer>
'@

# --- B4: e003 content (new text) - IsChecked -> IsEnabled (double space kept) ---
$textB4 = @'
<Bold>e003 Battle Board</Bold> <InlineUIContainer><Button Content='r2.12' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>
<LineBreak/><LineBreak/> <InlineUIContainer><CheckBox Name='SkipTutorial3' IsEnabled='True'  Content='Skip in Future' FontFamily='Courier New'  FontSize='12'></CheckBox></InlineUIContainer>
<LineBreak/><LineBreak/>
The battle board is an abstract display used to resolve engagements with enemy forces. Your tank is placed in the center of this display and the action of an engagement revolves around it through the use of pieces representing enemy units and other informational markers. A detailed explanation is given in <InlineUIContainer><Button Content='r5.0' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>. Click image to continue.
<LineBreak/><LineBreak/>
                                   <InlineUIContainer><Image Name='MapBattle'  Height='200' Width='200'></Image></InlineUIContainer>
'@

# --- B5: e004 content - IsChecked (double space) -> IsEnabled (single space) ---
$textB5 = @'
<Bold>e004 Tank Card</Bold> <InlineUIContainer><Button Content='r2.13' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>
<LineBreak/><LineBreak/> <InlineUIContainer><CheckBox Name='SkipTutorial4' IsEnabled='True' Content='Skip in Future' FontFamily='Courier New'  FontSize='12'></CheckBox></InlineUIContainer>
<LineBreak/><LineBreak/>
The upper right image is the Tank Card. The game starts with the basic M4 Sherman tank, i.e., Tank Card #1. 
The Tank Card shows the tank model and other important information regarding the tank. The use of the Tank Card is described in 
<InlineUIContainer><Button Content='r5.2' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>. Click image to continue.
<LineBreak/><LineBreak/>
                                 <InlineUIContainer><Image Name='m01'  Height='200' Width='200'></Image></InlineUIContainer>
'@

# --- B6: e005 content (new text) - CheckBox IsChecked attribute removed entirely ---
$textB6 = @'
<Bold>e005 After Action Report (AAR)</Bold> <InlineUIContainer><Button Content='r2.4' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>
<LineBreak/><LineBreak/> <InlineUIContainer><CheckBox Name='SkipTutorial5' Content='Skip in Future' FontFamily='Courier New'  FontSize='12'></CheckBox></InlineUIContainer>
<LineBreak/><LineBreak/>
The events of each engagement or day of battle are recorded as they unfold on the After Action Report. At this time, you may elect to change the name of the tank or the names of your crew by clicking on the appropriate location on the form. 
<LineBreak/><LineBreak/>When ready, click image below to assign crew ratings to your new crew per 
<InlineUIContainer><Button Content='r7.1' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>
<LineBreak/><LineBreak/>
                                            <InlineUIContainer><Image Name='Continue005' Height='100' Width='100'></Image></InlineUIContainer>
'@

$ws.Range("B1").Value = $textB1
$ws.Range("B2").Value = $textB2
$ws.Range("B3").Value = $textB3
$ws.Range("B4").Value = $textB4
$ws.Range("B5").Value = $textB5
$ws.Range("B6").Value = $textB6

# Row 6's content got shorter, so its height shrinks from 150 to 135.
$ws.Rows.Item(6).RowHeight = 135

# Scroll/selection state moved down to show the new row 4 (e003) and select B7.
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("B7").Select()
